$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused transaction rows (rows 4-7); rows 2:3 remain as the
# only transactions and everything below shifts up.
$ws.Range("A4:G7").EntireRow.Delete()

# --- Header row (row 1) ---
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "Time"
$ws.Range("D1").Value = "Purpose"
$ws.Range("E1").Value = "Trsct Type"
$ws.Range("F1").Value = "Trstcd amt"

# New "Balance" header in G1, copying the bold/centered/bordered header style
# from F1 so it matches the rest of the header row.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Balance"
$excel.CutCopyMode = $false

# --- Row 2 (first transaction) ---
$ws.Range("A2").Value = 1
# Assign the date as a formula-computed text, then paste it back as a plain
# value; this keeps it stored as a shared-string text constant ("20-June-2022")
# instead of Excel auto-converting the look-alike date text into a date serial.
$ws.Range("B2").Formula = '="20-June-2022"'
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("C2").Value = "10:41 AM"
$ws.Range("D2").Value = "dep"
$ws.Range("E2").Value = "Deposit [+]"
$ws.Range("F2").Value = 500
$ws.Range("G2").Value = 500

# --- Row 3 (second transaction) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Formula = '="20-June-2022"'
$ws.Range("B3").Copy()
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("C3").Value = "10:42 AM"
$ws.Range("D3").Value = "debit"
$ws.Range("E3").Value = "Withdrawn [-]"
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 400

$excel.CutCopyMode = $false
